$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture the current "Ledger Generation Date" text (2020-09-24) before it
#    is overwritten, so it can be reused for the new credit-note row below
#    (row 13) without Excel re-interpreting the literal as a date serial.
#    Stash it in row 1 (above every row we will shift) so a later row-insert
#    does not move it out from under us.
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy()
$ws.Range("G1").PasteSpecial(-4163)   # xlPasteValues -> stash as plain text

# ---------------------------------------------------------------------------
# 2. Update the "Ledger Generation Date" value to 2020-10-02, keeping it as
#    plain text with the original (default) cell style.
# ---------------------------------------------------------------------------
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2020-10-02"
$ws.Range("B5").Style = $ws.Range("A5").Style

# ---------------------------------------------------------------------------
# 3. Clear the "Particulars"/"Voucher" text for the two existing credit note
#    rows (11 and 12) - keep their formatting (style 2) intact.
# ---------------------------------------------------------------------------
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()

# ---------------------------------------------------------------------------
# 4. Insert a new row at position 13 (this pushes the old "TOTAL" row down to
#    row 14) and populate it as another credit note line, mirroring the
#    formatting of row 12.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# Copy formatting (borders/fill/font/number format) from row 12 onto row 13.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)   # xlPasteFormats

# A13: reuse the stashed "2020-09-24" text value (kept as plain text).
$ws.Range("G1").Copy()
$ws.Range("A13").PasteSpecial(-4163)       # xlPasteValues
$ws.Range("G1").ClearContents()

# B13 / C13 stay blank (numeric/empty), matching the cleared columns above.
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()

# D13: voucher type text.
$ws.Range("D13").Value = "CREDITNOTE"

# E13: credit amount (numeric 0).
$ws.Range("E13").Value = 0

# ---------------------------------------------------------------------------
# 5. Fix up the TOTAL row (now row 14) so its formula sums through the new
#    row 13.
# ---------------------------------------------------------------------------
$ws.Range("E14").Formula = "=SUM(E10:E13)"

Write-Output "edit applied"
